# "fix crear caja estimacion en presupuesto"
# The sheet used to list 5 generic "Estimacion N" line items; it now lists
# 3 named items (aaa, xd, Hola) with updated quantities/rates, and the
# summary box below (TOTAL / Linea Base de Costos / Presupuesto / ...) is
# recalculated to match the new totals, with new % cells (F column) added
# next to the "Reserva de gestion", "Ganancia" and "IGV" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 3 remaining line items (rows 2-4) ---
$ws.Range("A2").Value = "aaa"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 17904

$ws.Range("A3").Value = "xd"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 500
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1000

$ws.Range("A4").Value = "Hola"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1400
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 11200

# --- Remove the two extra "Estimacion" rows (old rows 5 and 6); this
#     shifts the summary box below up by two rows, preserving its labels ---
$ws.Rows("5:6").Delete()

# --- Recalculate the summary box values (now rows 6-14) ---
$ws.Range("E6").Value = 30104     # TOTAL
$ws.Range("E7").Value = 0         # Reserva de contingencia
$ws.Range("E8").Value = 30104     # Linea Base de Costos
$ws.Range("F9").Value = 0         # Reserva de gestion (%)
$ws.Range("E10").Value = 12500    # Presupuesto
$ws.Range("F11").Value = 0.05     # Ganancia (%)
$ws.Range("E12").Value = 42604    # Total con ganancia
$ws.Range("F13").Value = 0.15     # IGV (%)
$ws.Range("E14").Value = 42604.15 # Total

# --- Narrow column A (Partida) from 14 to 10 characters wide ---
$ws.Columns.Item(1).ColumnWidth = 9.140625
